$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 23
$prevRow = $row - 1

# Copy the date cell's formatting (built-in date/time number format) from the
# row above so the new row matches the existing style instead of minting a
# brand-new custom number format.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 42604.890231481484
$ws.Cells.Item($row, 2).Value = "Random"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 56
$ws.Cells.Item($row, 9).Value = 44
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 1
$ws.Cells.Item($row, 13).Value = 99
